$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 1.532141
$ws.Cells.Item(2, 8).Value2 = 4.596423
$ws.Cells.Item(2, 9).Value2 = 0.08900664250669833
$ws.Cells.Item(2, 10).Value2 = 0.08900664250669831
$ws.Cells.Item(2, 13).Value2 = 100.4511693333333
$ws.Cells.Item(2, 14).Value2 = 301.353508
$ws.Cells.Item(2, 15).Value2 = 0.6815338940941451
$ws.Cells.Item(2, 16).Value2 = 0.681533894094145
$ws.Cells.Item(2, 17).Value2 = 153.9053550335427
$ws.Cells.Item(2, 18).Value2 = 1385.148195301884
$ws.Cells.Item(2, 19).Value2 = 0.06066104366783558
$ws.Cells.Item(2, 20).Value2 = 0.06066104366783555
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 1.532141
$ws.Cells.Item(3, 8).Value2 = 4.596423
$ws.Cells.Item(3, 9).Value2 = 0.08900664250669833
$ws.Cells.Item(3, 10).Value2 = 0.08900664250669831
$ws.Cells.Item(3, 14).Value2 = 7.755446
$ws.Cells.Item(3, 15).Value2 = 0.01753953138921768
$ws.Cells.Item(3, 16).Value2 = 0.01753953138921768
$ws.Cells.Item(3, 17).Value2 = 3.960812263295333
$ws.Cells.Item(3, 18).Value2 = 35.647310369658
$ws.Cells.Item(3, 19).Value2 = 0.001561134800095112
$ws.Cells.Item(3, 20).Value2 = 0.001561134800095111
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 1.532141
$ws.Cells.Item(4, 8).Value2 = 4.596423
$ws.Cells.Item(4, 9).Value2 = 0.08900664250669833
$ws.Cells.Item(4, 10).Value2 = 0.08900664250669831
$ws.Cells.Item(4, 13).Value2 = 43.839503
$ws.Cells.Item(4, 14).Value2 = 131.518509
$ws.Cells.Item(4, 15).Value2 = 0.2974391178622877
$ws.Cells.Item(4, 16).Value2 = 0.2974391178622877
$ws.Cells.Item(4, 17).Value2 = 67.168299965923
$ws.Cells.Item(4, 18).Value2 = 604.514699693307
$ws.Cells.Item(4, 19).Value2 = 0.02647405723107635
$ws.Cells.Item(4, 20).Value2 = 0.02647405723107634
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 1.532141
$ws.Cells.Item(5, 8).Value2 = 4.596423
$ws.Cells.Item(5, 9).Value2 = 0.08900664250669833
$ws.Cells.Item(5, 10).Value2 = 0.08900664250669831
$ws.Cells.Item(5, 13).Value2 = 0.5140156666666666
$ws.Cells.Item(5, 14).Value2 = 1.542047
$ws.Cells.Item(5, 15).Value2 = 0.003487456654349595
$ws.Cells.Item(5, 16).Value2 = 0.003487456654349595
$ws.Cells.Item(5, 17).Value2 = 0.7875444775423333
$ws.Cells.Item(5, 18).Value2 = 7.087900297880999
$ws.Cells.Item(5, 19).Value2 = 0.0003104068076913007
$ws.Cells.Item(5, 20).Value2 = 0.0003104068076913005
$ws.Cells.Item(6, 9).Value2 = 0.6169137955113024
$ws.Cells.Item(6, 10).Value2 = 0.6169137955113023
$ws.Cells.Item(6, 13).Value2 = 100.4511693333333
$ws.Cells.Item(6, 14).Value2 = 301.353508
$ws.Cells.Item(6, 15).Value2 = 0.6815338940941451
$ws.Cells.Item(6, 16).Value2 = 0.681533894094145
$ws.Cells.Item(6, 17).Value2 = 1066.733156641787
$ws.Cells.Item(6, 18).Value2 = 9600.598409776081
$ws.Cells.Item(6, 19).Value2 = 0.4204476613752171
$ws.Cells.Item(6, 20).Value2 = 0.4204476613752169
$ws.Cells.Item(7, 9).Value2 = 0.6169137955113024
$ws.Cells.Item(7, 10).Value2 = 0.6169137955113023
$ws.Cells.Item(7, 14).Value2 = 7.755446
$ws.Cells.Item(7, 15).Value2 = 0.01753953138921768
$ws.Cells.Item(7, 16).Value2 = 0.01753953138921768
$ws.Cells.Item(7, 19).Value2 = 0.0108203788808119
$ws.Cells.Item(7, 20).Value2 = 0.0108203788808119
$ws.Cells.Item(8, 9).Value2 = 0.6169137955113024
$ws.Cells.Item(8, 10).Value2 = 0.6169137955113023
$ws.Cells.Item(8, 13).Value2 = 43.839503
$ws.Cells.Item(8, 14).Value2 = 131.518509
$ws.Cells.Item(8, 15).Value2 = 0.2974391178622877
$ws.Cells.Item(8, 16).Value2 = 0.2974391178622877
$ws.Cells.Item(8, 17).Value2 = 465.55009494826
$ws.Cells.Item(8, 18).Value2 = 4189.95085453434
$ws.Cells.Item(8, 19).Value2 = 0.1834942951339575
$ws.Cells.Item(8, 20).Value2 = 0.1834942951339575
$ws.Cells.Item(9, 9).Value2 = 0.6169137955113024
$ws.Cells.Item(9, 10).Value2 = 0.6169137955113023
$ws.Cells.Item(9, 13).Value2 = 0.5140156666666666
$ws.Cells.Item(9, 14).Value2 = 1.542047
$ws.Cells.Item(9, 15).Value2 = 0.003487456654349595
$ws.Cells.Item(9, 16).Value2 = 0.003487456654349595
$ws.Cells.Item(9, 17).Value2 = 5.458548250913333
$ws.Cells.Item(9, 18).Value2 = 49.12693425822
$ws.Cells.Item(9, 19).Value2 = 0.002151460121315957
$ws.Cells.Item(9, 20).Value2 = 0.002151460121315956
$ws.Cells.Item(10, 7).Value2 = 4.902263666666666
$ws.Cells.Item(10, 8).Value2 = 14.706791
$ws.Cells.Item(10, 9).Value2 = 0.2847871244569372
$ws.Cells.Item(10, 10).Value2 = 0.2847871244569371
$ws.Cells.Item(10, 13).Value2 = 100.4511693333333
$ws.Cells.Item(10, 14).Value2 = 301.353508
$ws.Cells.Item(10, 15).Value2 = 0.6815338940941451
$ws.Cells.Item(10, 16).Value2 = 0.681533894094145
$ws.Cells.Item(10, 17).Value2 = 492.4381176969809
$ws.Cells.Item(10, 18).Value2 = 4431.943059272829
$ws.Cells.Item(10, 19).Value2 = 0.1940920779190103
$ws.Cells.Item(10, 20).Value2 = 0.1940920779190103
$ws.Cells.Item(11, 7).Value2 = 4.902263666666666
$ws.Cells.Item(11, 8).Value2 = 14.706791
$ws.Cells.Item(11, 9).Value2 = 0.2847871244569372
$ws.Cells.Item(11, 10).Value2 = 0.2847871244569371
$ws.Cells.Item(11, 14).Value2 = 7.755446
$ws.Cells.Item(11, 15).Value2 = 0.01753953138921768
$ws.Cells.Item(11, 16).Value2 = 0.01753953138921768
$ws.Cells.Item(11, 17).Value2 = 12.67308038153178
$ws.Cells.Item(11, 18).Value2 = 114.057723433786
$ws.Cells.Item(11, 19).Value2 = 0.00499503270865749
$ws.Cells.Item(11, 20).Value2 = 0.004995032708657489
$ws.Cells.Item(12, 7).Value2 = 4.902263666666666
$ws.Cells.Item(12, 8).Value2 = 14.706791
$ws.Cells.Item(12, 9).Value2 = 0.2847871244569372
$ws.Cells.Item(12, 10).Value2 = 0.2847871244569371
$ws.Cells.Item(12, 13).Value2 = 43.839503
$ws.Cells.Item(12, 14).Value2 = 131.518509
$ws.Cells.Item(12, 15).Value2 = 0.2974391178622877
$ws.Cells.Item(12, 16).Value2 = 0.2974391178622877
$ws.Cells.Item(12, 17).Value2 = 214.9128027216243
$ws.Cells.Item(12, 18).Value2 = 1934.215224494619
$ws.Cells.Item(12, 19).Value2 = 0.08470683107700894
$ws.Cells.Item(12, 20).Value2 = 0.0847068310770089
$ws.Cells.Item(13, 7).Value2 = 4.902263666666666
$ws.Cells.Item(13, 8).Value2 = 14.706791
$ws.Cells.Item(13, 9).Value2 = 0.2847871244569372
$ws.Cells.Item(13, 10).Value2 = 0.2847871244569371
$ws.Cells.Item(13, 13).Value2 = 0.5140156666666666
$ws.Cells.Item(13, 14).Value2 = 1.542047
$ws.Cells.Item(13, 15).Value2 = 0.003487456654349595
$ws.Cells.Item(13, 16).Value2 = 0.003487456654349595
$ws.Cells.Item(13, 17).Value2 = 2.519840326797444
$ws.Cells.Item(13, 18).Value2 = 22.678562941177
$ws.Cells.Item(13, 19).Value2 = 0.0009931827522604319
$ws.Cells.Item(13, 20).Value2 = 0.0009931827522604315
$ws.Cells.Item(14, 5).Value2 = 2
$ws.Cells.Item(14, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(14, 7).Value2 = 0.159958
$ws.Cells.Item(14, 8).Value2 = 0.479874
$ws.Cells.Item(14, 9).Value2 = 0.009292437525062282
$ws.Cells.Item(14, 10).Value2 = 0.009292437525062281
$ws.Cells.Item(14, 13).Value2 = 100.4511693333333
$ws.Cells.Item(14, 14).Value2 = 301.353508
$ws.Cells.Item(14, 15).Value2 = 0.6815338940941451
$ws.Cells.Item(14, 16).Value2 = 0.681533894094145
$ws.Cells.Item(14, 17).Value2 = 16.06796814422133
$ws.Cells.Item(14, 18).Value2 = 144.611713297992
$ws.Cells.Item(14, 19).Value2 = 0.006333111132082257
$ws.Cells.Item(14, 20).Value2 = 0.006333111132082255
$ws.Cells.Item(15, 5).Value2 = 2
$ws.Cells.Item(15, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(15, 7).Value2 = 0.159958
$ws.Cells.Item(15, 8).Value2 = 0.479874
$ws.Cells.Item(15, 9).Value2 = 0.009292437525062282
$ws.Cells.Item(15, 10).Value2 = 0.009292437525062281
$ws.Cells.Item(15, 14).Value2 = 7.755446
$ws.Cells.Item(15, 15).Value2 = 0.01753953138921768
$ws.Cells.Item(15, 16).Value2 = 0.01753953138921768
$ws.Cells.Item(15, 17).Value2 = 0.4135152104226666
$ws.Cells.Item(15, 18).Value2 = 3.721636893804
$ws.Cells.Item(15, 19).Value2 = 0.0001629849996531741
$ws.Cells.Item(15, 20).Value2 = 0.0001629849996531741
$ws.Cells.Item(16, 5).Value2 = 2
$ws.Cells.Item(16, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(16, 7).Value2 = 0.159958
$ws.Cells.Item(16, 8).Value2 = 0.479874
$ws.Cells.Item(16, 9).Value2 = 0.009292437525062282
$ws.Cells.Item(16, 10).Value2 = 0.009292437525062281
$ws.Cells.Item(16, 13).Value2 = 43.839503
$ws.Cells.Item(16, 14).Value2 = 131.518509
$ws.Cells.Item(16, 15).Value2 = 0.2974391178622877
$ws.Cells.Item(16, 16).Value2 = 0.2974391178622877
$ws.Cells.Item(16, 17).Value2 = 7.012479220874
$ws.Cells.Item(16, 18).Value2 = 63.11231298786599
$ws.Cells.Item(16, 19).Value2 = 0.002763934420244945
$ws.Cells.Item(16, 20).Value2 = 0.002763934420244944
$ws.Cells.Item(17, 5).Value2 = 2
$ws.Cells.Item(17, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(17, 7).Value2 = 0.159958
$ws.Cells.Item(17, 8).Value2 = 0.479874
$ws.Cells.Item(17, 9).Value2 = 0.009292437525062282
$ws.Cells.Item(17, 10).Value2 = 0.009292437525062281
$ws.Cells.Item(17, 13).Value2 = 0.5140156666666666
$ws.Cells.Item(17, 14).Value2 = 1.542047
$ws.Cells.Item(17, 15).Value2 = 0.003487456654349595
$ws.Cells.Item(17, 16).Value2 = 0.003487456654349595
$ws.Cells.Item(17, 17).Value2 = 0.08222091800866665
$ws.Cells.Item(17, 18).Value2 = 0.7399882620779999
$ws.Cells.Item(17, 19).Value2 = 0.00003240697308190634
$ws.Cells.Item(17, 20).Value2 = 0.00003240697308190633
